$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 76 values/format down to rows 77 and 78, with incrementing dates.
for ($i = 77; $i -le 78; $i++) {
    $prev = $i - 1

    # Copy the previous row's formatting (border, font, number format) down,
    # then overwrite values: the date in column A increments by one day,
    # the rest of the row (B:J) repeats the previous row's values.
    $ws.Range("A$prev`:J$prev").Copy() | Out-Null
    $ws.Range("A$i`:J$i").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$i").Value2 = $ws.Range("A$prev").Value2 + 1
    $ws.Range("B$i").Value2 = $ws.Range("B$prev").Value2
    $ws.Range("C$i").Value2 = $ws.Range("C$prev").Value2
    $ws.Range("D$i").Value2 = $ws.Range("D$prev").Value2
    $ws.Range("E$i").Value2 = $ws.Range("E$prev").Value2
    $ws.Range("F$i").Value2 = $ws.Range("F$prev").Value2
    $ws.Range("G$i").Value2 = $ws.Range("G$prev").Value2
    $ws.Range("H$i").Value2 = $ws.Range("H$prev").Value2
    $ws.Range("I$i").Value2 = $ws.Range("I$prev").Value2
    $ws.Range("J$i").Value2 = $ws.Range("J$prev").Value2
}
